$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.839.62"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "3.167.42"
$ws.Range("E3").Value = "  -1.75%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.83"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.93"
$ws.Range("E6").Value = "  -4.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.593"
$ws.Range("E7").Value = "  -5.50%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -3.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.70"
$ws.Range("E10").Value = "  -1.22%  "
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").Value = "3.715.11"
$ws.Range("E12").Value = "  -2.10%  "
$ws.Range("E13").Value = "  -0.57%  "
$ws.Range("D14").Value = "64.745.32"
$ws.Range("E14").Value = "  -0.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.48"
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("D16").Value = "3.161.84"
$ws.Range("E16").Value = "  -2.63%  "
$ws.Range("E17").Value = "  -1.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "413.91"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.68"
$ws.Range("E19").Value = "  -1.06%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.29"
$ws.Range("E20").Value = "  -1.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.14"
$ws.Range("E21").Value = "  -0.43%  "
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.41"
$ws.Range("E23").Value = "  -2.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.202"
$ws.Range("E24").Value = "  -0.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.487"
$ws.Range("E25").Value = "  -1.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000104"
$ws.Range("E26").Value = "  -5.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.88"
$ws.Range("E27").Value = "  -1.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("E29").Value = "  -3.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "21.33"
$ws.Range("E30").Value = "  -1.75%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.93"
$ws.Range("E31").Value = "  -1.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.34"
$ws.Range("E32").Value = "  -0.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.14"
$ws.Range("E33").Value = "  -1.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "155.12"
$ws.Range("E34").Value = "  -0.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.36"
$ws.Range("E35").Value = "  -2.31%  "
$ws.Range("D36").Value = "2.718.63"
$ws.Range("E36").Value = "  -3.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.71"
$ws.Range("E37").Value = "  -1.92%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "23.77"
$ws.Range("E38").Value = "  -6.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.10"
$ws.Range("E39").Value = "  -2.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.706"
$ws.Range("E40").Value = "  -3.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0633"
$ws.Range("E41").Value = "  +1.08%  "
$ws.Range("E42").Value = "  -3.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0263"
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "291.56"
$ws.Range("E44").Value = "  -4.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.36"
$ws.Range("E45").Value = "  -3.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0984"
$ws.Range("E47").Value = "  -2.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.96"
$ws.Range("E48").Value = "  -10.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.48"
$ws.Range("E49").Value = "  +0.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.76"
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.898"
$ws.Range("E51").Value = "  -3.97%  "
